$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (FAPs->FAPs) to host the new
# FAPs->ECs entry, then strip the formatting that Insert() copied down
# from row 1 so the new row matches the plain (unstyled) data rows.
$ws.Rows("2:2").Insert()
$ws.Range("A2:T2").ClearFormats()

# New row 2: FAPs -> Tac2 -> Tacr1 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tac2"
$ws.Range("C2").Value = "Tacr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2004523333333333
$ws.Range("H2").Value = 0.601357
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2838346666666667
$ws.Range("N2").Value = 0.851504
$ws.Range("O2").Value = 0.7335492763611302
$ws.Range("P2").Value = 0.7335492763611302
$ws.Range("Q2").Value = 0.05689532121422223
$ws.Range("R2").Value = 0.512057890928
$ws.Range("S2").Value = 0.7335492763611302
$ws.Range("T2").Value = 0.7335492763611302

# Row 3 (previously row 2, FAPs->FAPs) gets refreshed values
$ws.Range("H3").Value = 0.601357
$ws.Range("M3").Value = 0.057966
$ws.Range("N3").Value = 0.173898
$ws.Range("O3").Value = 0.1498087525844245
$ws.Range("P3").Value = 0.1498087525844245
$ws.Range("Q3").Value = 0.011619419954
$ws.Range("R3").Value = 0.104574779586
$ws.Range("S3").Value = 0.1498087525844245
$ws.Range("T3").Value = 0.1498087525844245

# Row 4 (previously row 3, FAPs->MuSCs) gets refreshed values
$ws.Range("H4").Value = 0.601357
$ws.Range("O4").Value = 0.1166419710544452
$ws.Range("P4").Value = 0.1166419710544452
$ws.Range("Q4").Value = 0.009046948342888888
$ws.Range("R4").Value = 0.081422535086
$ws.Range("S4").Value = 0.1166419710544452
$ws.Range("T4").Value = 0.1166419710544452
